$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.915.52'
$ws.Range("E2").Value = '  +2.10%  '

$ws.Range("D3").Value = '1.813.43'
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = '  +0.75%  '

$ws.Range("D5").Value = "'311.67"
$ws.Range("E5").Value = '  +2.37%  '

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("D7").Value = "'0.4293"
$ws.Range("E7").Value = '  +0.51%  '

$ws.Range("D8").Value = "'0.3694"
$ws.Range("E8").Value = '  +1.98%  '

$ws.Range("D9").Value = "'0.07255"
$ws.Range("E9").Value = '  +3.01%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '2.143.11'
$ws.Range("E10").Value = '  +22.00%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = "'0.8665"
$ws.Range("E11").Value = '  +4.13%  '

$ws.Range("D12").Value = "'21.29"
$ws.Range("E12").Value = '  +5.56%  '

$ws.Range("D13").Value = "'6.639"
$ws.Range("E13").Value = '  +3.96%  '

$ws.Range("D14").Value = "'5.405"
$ws.Range("E14").Value = '  +3.31%  '

$ws.Range("D15").Value = "'0.06966"
$ws.Range("E15").Value = '  +2.60%  '

$ws.Range("D16").Value = "'80.88"
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").Value = "'0.000008899"
$ws.Range("E18").Value = '  +2.81%  '

$ws.Range("E19").Value = '  +0.41%  '

$ws.Range("D20").Value = "'15.28"
$ws.Range("E20").Value = '  +1.90%  '

$ws.Range("D21").Value = '26.969.93'
$ws.Range("E21").Value = '  +3.00%  '

$ws.Range("D22").Value = "'5.188"
$ws.Range("E22").Value = '  +3.72%  '

$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.381.45'
$ws.Range("E23").Value = '  +20.43%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = "'10.98"
$ws.Range("E24").Value = '  -1.22%  '

$ws.Range("D25").Value = "'154.30"
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("D26").Value = "'1.888"
$ws.Range("E26").Value = '  -0.76%  '

$ws.Range("E27").Value = '  +1.57%  '

$ws.Range("D28").Value = "'5.228"
$ws.Range("E28").Value = '  +4.16%  '

$ws.Range("D29").Value = "'1.942"
$ws.Range("E29").Value = '  +16.02%  '

$ws.Range("D30").Value = "'114.78"
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").Value = "'0.08958"
$ws.Range("E31").Value = '  +0.77%  '

$ws.Range("D32").Value = "'1.166"
$ws.Range("E32").Value = '  +5.04%  '

$ws.Range("D33").Value = "'0.7424"
$ws.Range("E33").Value = '  +2.95%  '

$ws.Range("D34").Value = "'4.432"
$ws.Range("E34").Value = '  +3.04%  '

$ws.Range("E35").Value = '  +2.24%  '

$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("D37").Value = "'1.121"
$ws.Range("E37").Value = '  +4.92%  '

$ws.Range("D38").Value = "'0.05233"
$ws.Range("E38").Value = '  +2.83%  '

$ws.Range("D39").Value = "'0.01925"
$ws.Range("E39").Value = '  +2.09%  '

$ws.Range("D40").Value = "'0.5099"
$ws.Range("E40").Value = '  +4.21%  '

$ws.Range("D41").Value = "'2.767"
$ws.Range("E41").Value = '  +10.92%  '

$ws.Range("D42").Value = "'0.1651"
$ws.Range("E42").Value = '  +3.16%  '

$ws.Range("D43").Value = "'6.474"
$ws.Range("E43").Value = '  +4.51%  '

$ws.Range("D44").Value = "'8.318"
$ws.Range("E44").Value = '  +4.33%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'10.49"
$ws.Range("E45").Value = '  +3.66%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = "'107.17"
$ws.Range("E46").Value = '  +2.18%  '

$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.4584"
$ws.Range("E48").Value = '  +2.88%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.06295"
$ws.Range("E49").Value = '  +1.67%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'1.648"
$ws.Range("E50").Value = '  +4.98%  '

$ws.Range("D51").Value = "'1.801"
$ws.Range("E51").Value = '  +3.64%  '
